$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$sh.Left = -107/12700
$sh.Top = -105/12700
$sh.Width = 9525212/12700
$sh.Height = 9525210/12700
Write-Host "Left $($sh.Left) Top $($sh.Top) Width $($sh.Width) Height $($sh.Height)"
